$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: force a Range to hold TEXT (inlineStr/shared-string) even when the
# value looks numeric (e.g. "23.45"), then strip the temporary "@" number
# format so the cell is left with the default style (no explicit style
# index), matching cells that were authored directly as inline strings.
# ---------------------------------------------------------------------------
function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# ---------------------------------------------------------------------------
# 1. "总计" (summary) sheet: insert a new row for 2022-Q3 at the top of the
#    data (row 2), pushing the existing quarters down by one row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()
# Re-use the formatting that's already on the row below (which used to be
# row 2, now shifted to row 3) so the new row matches the existing look
# (bold/centered/bordered "A" column cell, plain data cells).
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 6
$summary.Range("D2").Value = 1.8

# Column "A" is a plain 0-based row counter (0,1,2,...) independent of which
# quarter occupies the row - re-stamp it sequentially now that a row was
# inserted, rather than letting the insert carry the old counter values
# down with the shifted rows.
$summary.Range("A2").Value = 0
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5

# ---------------------------------------------------------------------------
# 2. Insert a brand-new "2022-Q3" worksheet right after "总计" holding the
#    per-fund breakdown for the new quarter (same layout as the other
#    quarterly sheets).
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q3"

# Header row (B1:H1) - text labels with the bold/centered/bordered style
# used by every other quarterly sheet.
Set-TextValue $newSheet.Range("B1") "基金代码"
Set-TextValue $newSheet.Range("C1") "基金名称"
Set-TextValue $newSheet.Range("D1") "基金规模"
Set-TextValue $newSheet.Range("E1") "股票总仓位"
Set-TextValue $newSheet.Range("F1") "仓位占比"
Set-TextValue $newSheet.Range("G1") "持有市值(亿元)"
Set-TextValue $newSheet.Range("H1") "仓位排名"

# Pull the header style off an existing quarterly sheet (its B1 cell already
# carries the correct bold/centered/bordered style) and stamp it across the
# whole new header row.
$existingQuarter = $wb.Worksheets.Item("2022-Q2")
$existingQuarter.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$fundRows = @(
    @(0, "519692", "交银成长混合A",         "23.45", "76.71", "3.15", "0.7387", 9),
    @(1, "519694", "交银蓝筹混合",           "16.10", "78.09", "3.08", "0.4959", 9),
    @(2, "011184", "东方阿尔法招阳混合A",     "6.44",  "89.08", "6.74", "0.4341", 7),
    @(3, "160919", "大成产业升级股票（LOF）",  "3.07",  "84.12", "3.92", "0.1203", 10),
    @(4, "011185", "东方阿尔法招阳混合C",     "0.10",  "89.08", "6.74", "0.0067", 7),
    @(5, "960016", "交银成长混合H",           "0.16",  "76.71", "3.15", "0.0050", 9)
)

$r = 2
foreach ($row in $fundRows) {
    $newSheet.Range("A$r").Value = $row[0]
    Set-TextValue $newSheet.Range("B$r") $row[1]
    Set-TextValue $newSheet.Range("C$r") $row[2]
    Set-TextValue $newSheet.Range("D$r") $row[3]
    Set-TextValue $newSheet.Range("E$r") $row[4]
    Set-TextValue $newSheet.Range("F$r") $row[5]
    Set-TextValue $newSheet.Range("G$r") $row[6]
    $newSheet.Range("H$r").Value = $row[7]
    $r++
}

# Column "A" index cells (A2:A7) use the same bold/centered/bordered style
# as the other quarterly sheets' "A" column - copy it across in one shot.
$existingQuarter.Range("A2").Copy()
$newSheet.Range("A2:A7").PasteSpecial(-4122)
